$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(82).Insert()

$ws.Range("A82").Value = 11
$ws.Range("B82").Value = "Vega Monumental Concepción"
$ws.Range("C82").Value = "Bíobío"
$ws.Range("D82").Value = 44985
$ws.Range("D82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103002
$ws.Range("J82").Value = "Ciruela"
$ws.Range("K82").Value = "Larry Ann"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 110
$ws.Range("N82").Value = 9000
$ws.Range("O82").Value = 10000
$ws.Range("P82").Value = 9455
$ws.Range("Q82").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R82").Value = "Región de O'Higgins"
$ws.Range("S82").Value = 525
$ws.Range("T82").Value = 18
